# Update symbol list: refresh Price / Volume(1h) values and fix the
# Bitrue/Mandala row ordering for rows 11-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/percentage columns hold numeric-looking text (the source sheet
# stores everything as text), so values are written with a leading
# apostrophe to force literal text and avoid Excel auto-converting them
# to numbers.
$textPrefix = "'"

# Row 2
$ws.Range("D2").Value = $textPrefix + '308.66'
$ws.Range("E2").Value = $textPrefix + '-2.87%'

# Row 3
$ws.Range("D3").Value = $textPrefix + '54.07'
$ws.Range("E3").Value = $textPrefix + '11.95%'

# Row 4
$ws.Range("D4").Value = $textPrefix + '5.093'
$ws.Range("E4").Value = $textPrefix + '-3.52%'

# Row 5
$ws.Range("D5").Value = $textPrefix + '0.07792'
$ws.Range("E5").Value = $textPrefix + '-1.95%'

# Row 6
$ws.Range("D6").Value = $textPrefix + '4.505'
$ws.Range("E6").Value = $textPrefix + '-1.94%'

# Row 7
$ws.Range("D7").Value = $textPrefix + '1.361'
$ws.Range("E7").Value = $textPrefix + '1.26%'

# Row 8
$ws.Range("D8").Value = $textPrefix + '1.571'
$ws.Range("E8").Value = $textPrefix + '-3.82%'

# Row 9
$ws.Range("D9").Value = $textPrefix + '0.1224'
$ws.Range("E9").Value = $textPrefix + '-4.45%'

# Row 10
$ws.Range("D10").Value = $textPrefix + '0.1993'
$ws.Range("E10").Value = $textPrefix + '2.10%'

# Row 11
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = $textPrefix + '0.04701'
$ws.Range("E11").Value = $textPrefix + '2.00%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = $textPrefix + '0.09375'
$ws.Range("E12").Value = $textPrefix + '-0.12%'

# Row 13
$ws.Range("D13").Value = $textPrefix + '0.1044'
$ws.Range("E13").Value = $textPrefix + '-0.36%'

# Row 14
$ws.Range("E14").Value = $textPrefix + '-4.53%'

# Row 15
$ws.Range("D15").Value = $textPrefix + '0.005765'
$ws.Range("E15").Value = $textPrefix + '-1.81%'

# Row 16
$ws.Range("E16").Value = $textPrefix + '2,012.71%'

# Row 17
$ws.Range("D17").Value = $textPrefix + '3.339'
$ws.Range("E17").Value = $textPrefix + '-0.01%'

# Row 18
$ws.Range("D18").Value = $textPrefix + '2.415'
$ws.Range("E18").Value = $textPrefix + '-0.34%'

# Row 19
$ws.Range("D19").Value = $textPrefix + '0.3443'
$ws.Range("E19").Value = $textPrefix + '-1.87%'

# Row 20
$ws.Range("D20").Value = $textPrefix + '8.048'
$ws.Range("E20").Value = $textPrefix + '-0.53%'

# Row 21
$ws.Range("D21").Value = $textPrefix + '0.1356'
$ws.Range("E21").Value = $textPrefix + '-2.01%'

# Row 22
$ws.Range("D22").Value = $textPrefix + '0.3074'
$ws.Range("E22").Value = $textPrefix + '-0.82%'

# Row 23
$ws.Range("D23").Value = $textPrefix + '0.04175'
$ws.Range("E23").Value = $textPrefix + '0.18%'

# Row 24
$ws.Range("D24").Value = $textPrefix + '0.001261'
$ws.Range("E24").Value = $textPrefix + '-4.48%'

# Row 25
$ws.Range("D25").Value = $textPrefix + '0.003959'
$ws.Range("E25").Value = $textPrefix + '-6.87%'

# Row 26
$ws.Range("D26").Value = $textPrefix + '0.0001353'
$ws.Range("E26").Value = $textPrefix + '0.07%'

# Row 38
$ws.Range("D38").Value = $textPrefix + '0.02625'
$ws.Range("E38").Value = $textPrefix + '-1.07%'

# Row 39
$ws.Range("D39").Value = $textPrefix + '0.05910'
$ws.Range("E39").Value = $textPrefix + '4.16%'

# Row 40
$ws.Range("D40").Value = $textPrefix + '0.01056'
$ws.Range("E40").Value = $textPrefix + '-1.97%'

# Row 41
$ws.Range("D41").Value = $textPrefix + '0.007913'
$ws.Range("E41").Value = $textPrefix + '-1.23%'

# Row 42
$ws.Range("D42").Value = $textPrefix + '0.1419'
$ws.Range("E42").Value = $textPrefix + '-1.06%'

# Row 43
$ws.Range("D43").Value = $textPrefix + '0.008193'
$ws.Range("E43").Value = $textPrefix + '10.01%'

# Row 44
$ws.Range("D44").Value = $textPrefix + '0.007839'
$ws.Range("E44").Value = $textPrefix + '-7.68%'

# Row 45
$ws.Range("D45").Value = $textPrefix + '0.3102'
$ws.Range("E45").Value = $textPrefix + '-1.93%'

# Row 46
$ws.Range("D46").Value = $textPrefix + '0.00007218'
$ws.Range("E46").Value = $textPrefix + '8.63%'

# Row 47
$ws.Range("D47").Value = $textPrefix + '0.00000000752'
$ws.Range("E47").Value = $textPrefix + '0.07%'

# Row 48
$ws.Range("D48").Value = $textPrefix + '0.05619'
$ws.Range("E48").Value = $textPrefix + '2.33%'

# Row 49
$ws.Range("D49").Value = $textPrefix + '0.002606'
$ws.Range("E49").Value = $textPrefix + '-34.92%'

# Row 50
$ws.Range("D50").Value = $textPrefix + '0.00002104'
$ws.Range("E50").Value = $textPrefix + '0.07%'

# Row 51
$ws.Range("D51").Value = $textPrefix + '0.0002004'
$ws.Range("E51").Value = $textPrefix + '0.07%'
